$d = $word.ActiveDocument

# Clear all existing paragraphs (delete from last to first to avoid offset shifts)
$paras = @()
foreach ($p in $d.Paragraphs) { $paras += $p }
for ($i = $paras.Count - 1; $i -ge 0; $i--) {
    $paras[$i].Range.Delete()
}

$xmlPayload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ContosoLearn Value Proposition</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ContosoLearn</w:t></w:r><w:r><w:t xml:space="preserve"> is not just an eLearning platform, it’s your </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>adaptive AI companion</w:t></w:r><w:r><w:t xml:space="preserve"> for personalized learning and skill development. Unlike other platforms, we offer a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>simplified yet powerful analytics system</w:t></w:r><w:r><w:t xml:space="preserve"> that provides actionable insights to enhance your learning experience without overwhelming you. Our platform is designed to </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>adapt to your learning style</w:t></w:r><w:r><w:t xml:space="preserve">, providing </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>customized recommendations</w:t></w:r><w:r><w:t xml:space="preserve"> based on your progress and preferences.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">We stand out in our commitment to </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>original content</w:t></w:r><w:r><w:t xml:space="preserve">. Our courses are not mere compilations of third-party </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>information, but</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> are </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>carefully curated and created by experts</w:t></w:r><w:r><w:t xml:space="preserve"> in the field, ensuring you receive unique and valuable knowledge. This commitment extends to providing a clear picture of your progress, areas of improvement, and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>next</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> steps, without getting lost in a sea of data.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In the face of stiff competition in the eLearning market, ContosoLearn stands out by continuously innovating and focusing on </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>bridging the gaps in traditional education</w:t></w:r><w:r><w:t>. We’re always looking for ways to make learning more efficient, effective, and enjoyable. Our focus isn’t just on what you learn, but also on how you learn.</w:t></w:r></w:p><w:p><w:r><w:t>In essence, ContosoLearn is more than an eLearning platform. It’s a personalized learning experience that adapts to you, grows with you, and empowers you to learn in a way that’s most effective for you. With ContosoLearn, you’re not just learning, you’re evolving. We are committed to helping you acquire new skills efficiently and effectively, making learning an enjoyable journey rather than a daunting task.</w:t></w:r></w:p><w:p/><w:p/><w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r = $d.Range(0, $d.Content.End)
$r.InsertXML($xmlPayload)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
Write-Host "Final text:" $d.Content.Text
